$wb = $excel.ActiveWorkbook

# Rename sheet "LEVEL 2" to "BoatsAttackGame"
$ws3 = $wb.Worksheets.Item("LEVEL 2")
$ws3.Name = "BoatsAttackGame"

# Update print area for the renamed sheet (keep same range, just refreshed with new sheet name)
$ws3.PageSetup.PrintArea = '$A$1:$BS$22'

# Update selections
$ws1 = $wb.Worksheets.Item("MAP")
$ws1.Activate()
$ws1.Range("E41").Select()

$ws3.Activate()
$ws3.Range("F45").Select()

# Update values (csv processing results)
$ws3.Range("BS3").Value = 21
$ws3.Range("B9").Value = 25
